$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.854.29'
$ws.Range("E2").Value = '  +4.38%  '
$ws.Range("D3").Value = '3.338.69'
$ws.Range("E3").Value = '  +4.33%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  +3.21%  '
$ws.Range("D6").Value = '''151.92'
$ws.Range("E6").Value = '  +4.48%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +2.10%  '
$ws.Range("E9").Value = '  +2.17%  '
$ws.Range("E10").Value = '  +3.74%  '
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("D12").Value = '3.915.55'
$ws.Range("E12").Value = '  +4.45%  '
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").Value = '''26.81'
$ws.Range("E15").Value = '  +2.43%  '
$ws.Range("D16").Value = '62.903.90'
$ws.Range("E16").Value = '  +4.55%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '''6.48'
$ws.Range("E17").Value = '  +4.32%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.271.84'
$ws.Range("E18").Value = '  +1.82%  '
$ws.Range("D19").Value = '''13.73'
$ws.Range("E19").Value = '  +4.57%  '
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("D21").Value = '''387.80'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("E23").Value = '  +1.39%  '
$ws.Range("D24").Value = '''70.74'
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("D25").Value = '''0.182'
$ws.Range("E25").Value = '  +5.04%  '
$ws.Range("D26").Value = '''8.80'
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").Value = '0.0₃0957'
$ws.Range("E27").Value = '  +5.78%  '
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''1.98'
$ws.Range("E29").Value = '  +3.41%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '''6.45'
$ws.Range("E30").Value = '  +4.14%  '
$ws.Range("D31").Value = '''22.97'
$ws.Range("E31").Value = '  +2.26%  '
$ws.Range("E32").Value = '  +2.34%  '
$ws.Range("E33").Value = '  +5.25%  '
$ws.Range("E34").Value = '  +2.22%  '
$ws.Range("D35").Value = '''160.68'
$ws.Range("E35").Value = '  +2.53%  '
$ws.Range("E36").Value = '  +9.25%  '
$ws.Range("E37").Value = '  +10.98%  '
$ws.Range("D38").Value = '''27.17'
$ws.Range("E38").Value = '  +5.58%  '
$ws.Range("D39").Value = '2.830.07'
$ws.Range("E39").Value = '  +2.33%  '
$ws.Range("D40").Value = '''0.0734'
$ws.Range("E40").Value = '  +2.88%  '
$ws.Range("D41").Value = '''0.0310'
$ws.Range("E41").Value = '  +7.89%  '
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("E43").Value = '  +2.49%  '
$ws.Range("E44").Value = '  +2.00%  '
$ws.Range("D45").Value = '''1.04'
$ws.Range("E45").Value = '  +2.60%  '
$ws.Range("D46").Value = '3.384.98'
$ws.Range("E46").Value = '  +4.48%  '
$ws.Range("D47").Value = '''21.91'
$ws.Range("E47").Value = '  +6.14%  '
$ws.Range("E48").Value = '  +3.20%  '
$ws.Range("D49").Value = '''6.27'
$ws.Range("E49").Value = '  +1.14%  '
$ws.Range("D50").Value = '''0.801'
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("D51").Value = '''282.60'
$ws.Range("E51").Value = '  +7.25%  '
